$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.726.43'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.895.17'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.86'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4918'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2939'
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06752'
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.894.45'
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.18'
$ws.Range('E11').Value = '  +4.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07246'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '90.87'
$ws.Range('E13').Value = '  +5.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6765'
$ws.Range('E14').Value = '  +0.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.045'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.681.57'
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008016'
$ws.Range('E17').Value = '  +2.18%  '
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.15'
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.139.45'
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.811'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '191.42'
$ws.Range('E23').Value = '  +32.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.107'
$ws.Range('E24').Value = '  +3.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.403'
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.02'
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.94'
$ws.Range('E27').Value = '  +11.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.894'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.413'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.307'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09074'
$ws.Range('E31').Value = '  +2.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.008'
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05300'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7470'
$ws.Range('E34').Value = '  +2.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.103'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.754'
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01831'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.686'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9350'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.116'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4413'
$ws.Range('E41').Value = '  +3.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.32'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.737'
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1359'
$ws.Range('E45').Value = '  +5.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.520'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05875'
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('E48').Value = '  +4.67%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3965'
$ws.Range('E49').Value = '  +4.59%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.94'
$ws.Range('E50').Value = '  +3.07%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.427'
$ws.Range('E51').Value = '  +5.29%  '
